# Auto-generated edit script applying the Diabolos_Profits.xlsx market-data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 210.55556
$ws.Range("J58").Value = 500
$ws.Range("L58").Value = 1500
$ws.Range("N58").Value = -1800
$ws.Range("H62").Value = 31226
$ws.Range("I62").Value = 3370.1428
$ws.Range("J62").Value = 47475.25
$ws.Range("K62").Value = 3370.1428
$ws.Range("L62").Value = 47475.25
$ws.Range("M62").Value = -2746.1428
$ws.Range("N62").Value = -48723.25
$ws.Range("H65").Value = 31226
$ws.Range("I65").Value = 3370.1428
$ws.Range("J65").Value = 47475.25
$ws.Range("K65").Value = 16850.714
$ws.Range("L65").Value = 237376.25
$ws.Range("M65").Value = -13730.714
$ws.Range("N65").Value = -243616.25
$ws.Range("H76").Value = 3377841
$ws.Range("I76").Value = 4077.8
$ws.Range("J76").Value = 5627016.5
$ws.Range("K76").Value = 4077.8
$ws.Range("L76").Value = 5627016.5
$ws.Range("M76").Value = -3762.8
$ws.Range("N76").Value = -5627646.5
$ws.Range("H79").Value = 3377841
$ws.Range("I79").Value = 4077.8
$ws.Range("J79").Value = 5627016.5
$ws.Range("K79").Value = 4077.8
$ws.Range("L79").Value = 5627016.5
$ws.Range("M79").Value = -2985.8
$ws.Range("N79").Value = -5629200.5
$ws.Range("H86").Value = 11170377
$ws.Range("I86").Value = 2472.625
$ws.Range("J86").Value = 20104700
$ws.Range("K86").Value = 2472.625
$ws.Range("L86").Value = 20104700
$ws.Range("M86").Value = -1349.625
$ws.Range("N86").Value = -20106946
$ws.Range("H88").Value = 2144.3635
$ws.Range("J88").Value = 2648.5
$ws.Range("L88").Value = 2648.5
$ws.Range("N88").Value = -3460.5
$ws.Range("H89").Value = 11170377
$ws.Range("I89").Value = 2472.625
$ws.Range("J89").Value = 20104700
$ws.Range("K89").Value = 12363.125
$ws.Range("L89").Value = 100523500
$ws.Range("M89").Value = -6747.125
$ws.Range("N89").Value = -100534732
$ws.Range("H91").Value = 2144.3635
$ws.Range("J91").Value = 2648.5
$ws.Range("L91").Value = 2648.5
$ws.Range("N91").Value = -5456.5
$ws.Range("H106").Value = 55494.895
$ws.Range("I106").Value = 2976.8823
$ws.Range("J106").Value = 501898
$ws.Range("K106").Value = 2976.8823
$ws.Range("L106").Value = 501898
$ws.Range("M106").Value = -2345.8823
$ws.Range("N106").Value = -503160
$ws.Range("H116").Value = 41740976
$ws.Range("I116").Value = 27890990
$ws.Range("J116").Value = 66670948
$ws.Range("K116").Value = 27890990
$ws.Range("L116").Value = 66670948
$ws.Range("M116").Value = -27887548
$ws.Range("N116").Value = -66677832
$ws.Range("H137").Value = 5309.5557
$ws.Range("I137").Value = 3239.25
$ws.Range("J137").Value = 6965.8
$ws.Range("K137").Value = 9717.75
$ws.Range("L137").Value = 20897.4
$ws.Range("M137").Value = -7167.75
$ws.Range("N137").Value = -25997.4
$ws.Range("H138").Value = 4368.643
$ws.Range("I138").Value = 4924.385
$ws.Range("J138").Value = 4241.8945
$ws.Range("K138").Value = 14773.155
$ws.Range("L138").Value = 12725.6835
$ws.Range("M138").Value = -9633.155000000001
$ws.Range("N138").Value = -23005.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15491.191
$ws.Range("I32").Value = 10727.968
$ws.Range("J32").Value = 45499.5
$ws.Range("K32").Value = 10727.968
$ws.Range("L32").Value = 45499.5
$ws.Range("M32").Value = -10440.968
$ws.Range("N32").Value = -46073.5
$ws.Range("H52").Value = 99166.664
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 99166.664
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 99166.664
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = -99802.664
$ws.Range("H61").Value = 419887.72
$ws.Range("I61").Value = 2696.6875
$ws.Range("J61").Value = 1254269.8
$ws.Range("K61").Value = 2696.6875
$ws.Range("L61").Value = 1254269.8
$ws.Range("M61").Value = -2484.6875
$ws.Range("N61").Value = -1254693.8
$ws.Range("H74").Value = 4358.381
$ws.Range("I74").Value = 3595.1428
$ws.Range("K74").Value = 3595.1428
$ws.Range("M74").Value = -2721.1428
$ws.Range("H77").Value = 4358.381
$ws.Range("I77").Value = 3595.1428
$ws.Range("K77").Value = 17975.714
$ws.Range("M77").Value = -13607.714
$ws.Range("H132").Value = 872816.0600000001
$ws.Range("I132").Value = 558183.2
$ws.Range("J132").Value = 2005494.6
$ws.Range("K132").Value = 1674549.6
$ws.Range("L132").Value = 6016483.800000001
$ws.Range("M132").Value = -1672019.6
$ws.Range("N132").Value = -6021543.800000001
$ws.Range("H136").Value = 419887.72
$ws.Range("I136").Value = 2696.6875
$ws.Range("J136").Value = 1254269.8
$ws.Range("K136").Value = 8090.0625
$ws.Range("L136").Value = 3762809.4
$ws.Range("M136").Value = -5540.0625
$ws.Range("N136").Value = -3767909.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 48999.668
$ws.Range("J35").Value = 48999.668
$ws.Range("L35").Value = 48999.668
$ws.Range("N35").Value = -49619.668
$ws.Range("H134").Value = 3222.697
$ws.Range("I134").Value = 2736.1538
$ws.Range("J134").Value = 5029.857
$ws.Range("K134").Value = 8208.4614
$ws.Range("L134").Value = 15089.571
$ws.Range("M134").Value = -5673.4614
$ws.Range("N134").Value = -20159.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2900.2075
$ws.Range("I31").Value = 2269.1
$ws.Range("J31").Value = 4842.077
$ws.Range("K31").Value = 2269.1
$ws.Range("L31").Value = 4842.077
$ws.Range("M31").Value = -1974.1
$ws.Range("N31").Value = -5432.077
$ws.Range("H34").Value = 2900.2075
$ws.Range("I34").Value = 2269.1
$ws.Range("J34").Value = 4842.077
$ws.Range("K34").Value = 2269.1
$ws.Range("L34").Value = 4842.077
$ws.Range("M34").Value = -2067.1
$ws.Range("N34").Value = -5246.077
$ws.Range("H58").Value = 235541.8
$ws.Range("I58").Value = 1645.0416
$ws.Range("J58").Value = 530990.3
$ws.Range("K58").Value = 1645.0416
$ws.Range("L58").Value = 530990.3
$ws.Range("M58").Value = -1442.0416
$ws.Range("N58").Value = -531396.3
$ws.Range("H132").Value = 235069.75
$ws.Range("I132").Value = 2513.5
$ws.Range("K132").Value = 7540.5
$ws.Range("M132").Value = -5010.5
$ws.Range("H134").Value = 5268.5386
$ws.Range("I134").Value = 4434.3105
$ws.Range("J134").Value = 7687.8
$ws.Range("K134").Value = 13302.9315
$ws.Range("L134").Value = 23063.4
$ws.Range("M134").Value = -10767.9315
$ws.Range("N134").Value = -28133.4
$ws.Range("H136").Value = 235541.8
$ws.Range("I136").Value = 1645.0416
$ws.Range("J136").Value = 530990.3
$ws.Range("K136").Value = 4935.1248
$ws.Range("L136").Value = 1592970.9
$ws.Range("M136").Value = -2385.1248
$ws.Range("N136").Value = -1598070.9
$ws.Range("H141").Value = 208291.92
$ws.Range("J141").Value = 214887.42
$ws.Range("L141").Value = 214887.42
$ws.Range("N141").Value = -225247.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 128.25
$ws.Range("I23").Value = 64
$ws.Range("K23").Value = 192
$ws.Range("M23").Value = 43
$ws.Range("H113").Value = 2611.875
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2611.875
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7835.625
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -12175.625
$ws.Range("H122").Value = 909.9
$ws.Range("J122").Value = 997.75
$ws.Range("L122").Value = 8979.75
$ws.Range("N122").Value = -13879.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 218493.33
$ws.Range("J39").Value = 218493.33
$ws.Range("L39").Value = 218493.33
$ws.Range("N39").Value = -219557.33
$ws.Range("H80").Value = 2567.4546
$ws.Range("J80").Value = 2964.5
$ws.Range("L80").Value = 2964.5
$ws.Range("N80").Value = -4960.5
$ws.Range("H83").Value = 2567.4546
$ws.Range("J83").Value = 2964.5
$ws.Range("L83").Value = 14822.5
$ws.Range("N83").Value = -24806.5
$ws.Range("H132").Value = 275602.97
$ws.Range("I132").Value = 316916.1
$ws.Range("J132").Value = 11199
$ws.Range("K132").Value = 950748.2999999999
$ws.Range("L132").Value = 33597
$ws.Range("M132").Value = -948218.2999999999
$ws.Range("N132").Value = -38657

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 931.6667
$ws.Range("I22").Value = 789.5
$ws.Range("K22").Value = 789.5
$ws.Range("M22").Value = -494.5
$ws.Range("H27").Value = 931.6667
$ws.Range("I27").Value = 789.5
$ws.Range("K27").Value = 789.5
$ws.Range("M27").Value = -682.5
$ws.Range("H122").Value = 7853.75
$ws.Range("I122").Value = 7853
$ws.Range("K122").Value = 23559
$ws.Range("M122").Value = -21109
$ws.Range("H132").Value = 131281.53
$ws.Range("I132").Value = 220683.67
$ws.Range("J132").Value = 6660.364
$ws.Range("K132").Value = 662051.01
$ws.Range("L132").Value = 19981.092
$ws.Range("M132").Value = -659521.01
$ws.Range("N132").Value = -25041.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 147777
$ws.Range("I64").Value = 147777
$ws.Range("K64").Value = 147777
$ws.Range("M64").Value = -147529
$ws.Range("H67").Value = 147777
$ws.Range("I67").Value = 147777
$ws.Range("K67").Value = 147777
$ws.Range("M67").Value = -146919
$ws.Range("H93").Value = 63999
$ws.Range("I93").Value = 63999
$ws.Range("K93").Value = 63999
$ws.Range("M93").Value = -61503
$ws.Range("H132").Value = 339799.1
$ws.Range("I132").Value = 406628.12
$ws.Range("J132").Value = 5654
$ws.Range("K132").Value = 1219884.36
$ws.Range("L132").Value = 16962
$ws.Range("M132").Value = -1217354.36
$ws.Range("N132").Value = -22022
$ws.Range("H135").Value = 45357
$ws.Range("J135").Value = 45357
$ws.Range("L135").Value = 45357
$ws.Range("N135").Value = -55497
$ws.Range("H136").Value = 7795.2104
$ws.Range("I136").Value = 7700.645
$ws.Range("J136").Value = 8214
$ws.Range("K136").Value = 23101.935
$ws.Range("L136").Value = 24642
$ws.Range("M136").Value = -20551.935
$ws.Range("N136").Value = -29742
